$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to refreshed values.
# Cells whose new text would otherwise be auto-parsed as a number
# are forced to Text format first so they stay as literal strings
# (matching the source data which stores every Price/Volume cell as text).

$ws.Range('D2').Value = '65.648.23'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '2.670.93'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.92'
$ws.Range('E5').Value = '  -1.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.39'
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.613'
$ws.Range('E8').Value = '  +4.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.128'
$ws.Range('E9').Value = '  +3.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.86'
$ws.Range('E10').Value = '  -2.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.398'
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.16'
$ws.Range('E13').Value = '  -4.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000195'
$ws.Range('E14').Value = '  -4.16%  '
$ws.Range('D15').Value = '3.153.57'
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('D16').Value = '65.523.29'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').Value = '2.680.76'
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.84'
$ws.Range('E18').Value = '  +0.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.79'
$ws.Range('E19').Value = '  -2.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.56'
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '352.38'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.69'
$ws.Range('E23').Value = '  -1.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000111'
$ws.Range('E24').Value = '  +4.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.61'
$ws.Range('E25').Value = '  -3.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.68'
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('E27').Value = '  -2.51%  '
$ws.Range('E28').Value = '  -5.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.02'
$ws.Range('E29').Value = '  -5.74%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('E31').Value = '  -2.88%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '527.75'
$ws.Range('E32').Value = '  -3.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.76'
$ws.Range('E33').Value = '  -3.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.53'
$ws.Range('E34').Value = '  +1.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.48'
$ws.Range('E35').Value = '  -3.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.422'
$ws.Range('E36').Value = '  -2.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.53'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '157.77'
$ws.Range('E39').Value = '  -3.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.94'
$ws.Range('E40').Value = '  -2.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '163.21'
$ws.Range('E42').Value = '  -5.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.12'
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('E44').Value = '  +2.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0608'
$ws.Range('E45').Value = '  -1.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.70'
$ws.Range('E46').Value = '  -4.04%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0264'
$ws.Range('E47').Value = '  +15.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.638'
$ws.Range('E48').Value = '  -2.52%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0257'
$ws.Range('E49').Value = '  -3.45%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0995'
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.09'
$ws.Range('E51').Value = '  -4.74%  '
